$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Infinity"/"-Infinity" text cells to use the compact "∞"/"-∞"
# notation, matching the new Numbers.toString() formatting for double/float
# infinities (EPBDS-12566). A leading apostrophe preserves the text
# quote-prefix formatting these cells already had.
$ws.Range("L12").Value = "'∞"
$ws.Range("M12").Value = "'-∞"
$ws.Range("L24").Value = "'∞"
$ws.Range("M24").Value = "'-∞"

# Update the active selection on the sheet
$null = $ws.Range("M23").Select()
